# Append two new Lancers job listings (scraped 2025-09-27 01:13:25 JST)
# to the top "ランサーズ" sheet, ahead of the lowest-priority existing
# row, and refresh the "取得日時" (fetched-at) timestamp on every
# existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-09-27 01:13:25"

# 1) Refresh the fetched-at timestamp (column A) on every existing data
#    row (2-16) before the new rows shuffle row 16 down to row 18.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# 2) Make room for two new rows right above the current last row (16),
#    pushing the current row 16 ("スーパードルフィー...") down to row 18.
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()

# 3) Row 16: Android kotlin job listing.
$ws.Cells.Item(16, 1).Value = $newTimestamp
$ws.Cells.Item(16, 2).Value = "Android kotlinの画像ファイル 拡張子取得"
$ws.Cells.Item(16, 3).Value = "システム開発"
$ws.Cells.Item(16, 4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(16, 5).Value = "期限情報なし"
$ws.Cells.Item(16, 6).Value = "https://www.lancers.jp/work/detail/5401572"
$ws.Cells.Item(16, 7).Value = 10

# 4) Row 17: WordPress local-environment support job listing.
$ws.Cells.Item(17, 1).Value = $newTimestamp
$ws.Cells.Item(17, 2).Value = "【急募】ワードプレスのローカル環境構築をサポートしてくれる方"
$ws.Cells.Item(17, 3).Value = "システム開発"
$ws.Cells.Item(17, 4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(17, 5).Value = "期限情報なし"
$ws.Cells.Item(17, 6).Value = "https://www.lancers.jp/work/detail/5401534"
$ws.Cells.Item(17, 7).Value = 10

# 5) The row pushed down to 18 lost its hyperlink object (the insert
#    only relocates cell content, not the hyperlink anchor), and the
#    two freshly-populated URL cells above have none yet - add both
#    now, in ascending row order, so relationship ids come out
#    rId16 (F17), rId17 (F18), matching the order new links are
#    appended to the part.
$ws.Hyperlinks.Add($ws.Cells.Item(17, 6), "https://www.lancers.jp/work/detail/5401534")
$ws.Hyperlinks.Add($ws.Cells.Item(18, 6), "https://www.lancers.jp/work/detail/5400988")

# Hyperlinks.Add() stamps a fresh ad-hoc "applyFont" style variant; snap
# both cells back onto the workbook's single shared Hyperlink style (the
# same index F2:F16 already use) instead of leaving a duplicate behind.
$ws.Cells.Item(17, 6).Style = "Hyperlink"
$ws.Cells.Item(18, 6).Style = "Hyperlink"
